$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player/position rows (target rows 604-691), listed in the order the
# author originally entered them (name/position pairs keep the row they
# ended up on after the sheet was later resorted).
$newPlayers = @(
    @(604,'Braxton Key','SG'),
    @(605,'Eugene Omoruyi','SF'),
    @(606,'Trayvon Palmer','SG'),
    @(607,'Jose Alvarado','PG'),
    @(608,'Ade Murkey','SG'),
    @(609,'Scottie Lewis','SG'),
    @(610,'Davon Reed','SG'),
    @(611,'Ish Wainright','SG'),
    @(612,'Duane Washington Jr.','PG'),
    @(613,'Javonte Smart','PG'),
    @(614,'Matt Mooney','SG'),
    @(615,'Arnoldas Kulboka','SF'),
    @(616,'Micah Potter','PF'),
    @(617,'Marko Simonovic','C'),
    @(618,'Jon Teske','C'),
    @(619,'Haywood Highsmith','SF'),
    @(620,'Cheick Diallo','PF'),
    @(621,'Rob Edwards','SG'),
    @(622,'Olivier Sarr','SG'),
    @(623,'Moses Wright','PF'),
    @(624,'Darren Collison','PG'),
    @(625,'Daishen Nix','SG'),
    @(626,'Isaiah Livers','SF'),
    @(627,'Justin Anderson','SF'),
    @(628,'Paris Bass','SF'),
    @(691,'Craig Sword','SG'),
    @(690,'JaQuori McLaughlin','PG'),
    @(689,'Klay Thompson','SG'),
    @(688,'Cade Cunningham','SG'),
    @(687,'Shaq Buchanan','SG'),
    @(686,'Yves Pons','SF'),
    @(685,'Joe Wieskamp','SG'),
    @(684,'Xavier Moon','SG'),
    @(683,'Jaysean Paige','PG'),
    @(682,'George King','SF'),
    @(681,'Ahmad Caver','PG'),
    @(680,'Jay Huff','C'),
    @(679,'Lance Stephenson','SG'),
    @(678,'Keifer Sykes','PG'),
    @(677,'Hassani Gravett','PG'),
    @(676,'BJ Johnson','SF'),
    @(675,'Brandon Knight','PG'),
    @(674,'Cameron McGriff','SF'),
    @(673,'Trevon Scott','SF'),
    @(672,'MJ Walker','SG'),
    @(671,'Jarron Cumberland','SG'),
    @(670,'Lindell Wigginton','PG'),
    @(669,'Jeff Dowtin','PG'),
    @(668,'Jordan Schakel','SG'),
    @(667,'Jemerrio Jones','SF'),
    @(666,'CJ Miles','SF'),
    @(661,'Neemias Queta','C'),
    @(663,'Aaron Wiggins','SG'),
    @(662,'Enes Freedom','C'),
    @(664,'Trevelin Queen','SG'),
    @(665,'Greg Monroe','C'),
    @(660,'Aleem Ford','SF'),
    @(659,'Carlik Jones','PG'),
    @(658,'Admiral Schofield','SF'),
    @(657,'Tyler Hall','SG'),
    @(656,'David Duke Jr.','SF'),
    @(655,'David Johnson','SG'),
    @(654,'Terry Taylor','SG'),
    @(653,'Mac McClung','PG'),
    @(652,'Jonathan Kuminga','SF'),
    @(629,'Xavier Sneed','SF'),
    @(630,'Scotty Hopson','SF'),
    @(631,'Joe Johnson','SG'),
    @(632,'Zylan Cheatham','SF'),
    @(633,'Kessler Edwards','SF'),
    @(634,'Trendon Watford','SF'),
    @(635,'Malcolm Hill','SG'),
    @(636,'Brandon Williams','PG'),
    @(637,'Keon Johnson','SF'),
    @(638,'Emmanuel Mudiay','PG'),
    @(639,'Emanuel Terry','C'),
    @(640,'Jaylen Morris','SG'),
    @(641,'Myles Powell','SG'),
    @(642,'Derrick Walton Jr.','PG'),
    @(643,'Jordan Goodwin','SG'),
    @(644,'Malik Newman','SG'),
    @(645,'Chaundee Brown Jr.','SG'),
    @(646,'Sam Hauser','PF'),
    @(647,'Cat Barber','PG'),
    @(648,'DeJon Jarreau','PG'),
    @(649,'Javin DeLaurier','C'),
    @(650,'Jaime Echenique','C'),
    @(651,'Nik Stauskas','SG')
)

foreach ($entry in $newPlayers) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
}

$lastRow = 691
$nextRow = $lastRow + 1

# Move the selection to the first empty row below the new data, and scroll
# the sheet so that region is in view (mirrors the author's saved view).
$ws.Range("B" + $nextRow).Select() | Out-Null
$excel.ActiveWindow.ScrollRow = $lastRow - 20
$excel.ActiveWindow.ScrollColumn = 1

# The workbook was re-printed in portrait orientation.
$ws.PageSetup.Orientation = 1
